# Modified the resource for Asset services
#
# - paymentMethod!C3 and paymentMethod!C4 change from the numeric
#   placeholder 102 to the text "state bank of india".
# - The active sheet/tab moves from "approvalDetails" to "paymentMethod",
#   with the selection on paymentMethod landing on C4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("paymentMethod")

$ws.Range("C3").Value = "state bank of india"
$ws.Range("C4").Value = "state bank of india"

# Make paymentMethod the active/selected sheet (mirrors activeTab moving
# from approvalDetails to paymentMethod in the workbook view), and leave
# the selection on C4.
$ws.Activate()
$ws.Range("C4").Select()
